$p = $ppt.ActivePresentation

# 1) Change the table style id on slide 16's table (3rd shape) to the built-in style guid.
$s = $p.Slides.Item(16)
$tblShape = $s.Shapes.Item(3)
$tblShape.Table.ApplyStyle("{379CD035-E40D-4B0B-BC47-F42C731B15D3}")
